# Swap the betting-data payload (columns B:AD) between pairs of adjacent
# rows. The leading rank/index in column A stays put; everything else
# (match id, teams, result, odds, ...) trades places between the two rows
# of each pair.
#
# NB: multi-cell Range.Value reads come back as an opaque placeholder in
# this host, so values are swapped column-by-column using single-cell
# Value2 reads (which are correctly typed) paired with single-cell Value
# writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(24, 73, 109, 204, 219)
$columns = @(2..30)   # B (2) through AD (30)

foreach ($r1 in $rowPairs) {
    $r2 = $r1 + 1

    foreach ($col in $columns) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}
